# Apply the commit's changes to the presentation using the PowerPoint COM object model.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: update the presenter/date line on the title slide.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
foreach ($shp in $s1.Shapes) {
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*Presenter Name*") {
            $tr.Text = $tr.Text.Replace("November 22, 2025", "November 24, 2025")
        }
    }
}

# ---------------------------------------------------------------------------
# Slide 9: Investment Summary table - remove Professional Services costs,
# rename Software -> Software Licenses, and fix the TOTAL row.
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
foreach ($shp in $s9.Shapes) {
    if ($shp.HasTable) {
        $tbl = $shp.Table

        # Row 2 = "Professional Services" -> Year 1 List / Year 1 Net / 3-Year Total become $0
        $tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "`$0"
        $tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = "`$0"
        $tbl.Cell(2, 7).Shape.TextFrame.TextRange.Text = "`$0"

        # Row 4 = "Software" -> "Software Licenses"
        $tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Software Licenses"

        # Row 6 = "TOTAL" -> update Year 1 List / Year 1 Net / 3-Year Total
        $tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "`$310,200"
        $tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "`$281,200"
        $tbl.Cell(6, 7).Shape.TextFrame.TextRange.Text = "`$438,400"
    }
}
